# Append three new monthly data rows (225-227) to the Romania M2 sheet,
# mirroring the existing "ECONOMICS:ROM2" rows (same layout/format as row 224).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 225; Date = 45108.41666666666; Value = 622009300000 },
    @{ Row = 226; Date = 45139.41666666666; Value = 626029700000 },
    @{ Row = 227; Date = 45170.41666666666; Value = 640762400000 }
)

foreach ($item in $newRows) {
    $r = $item.Row

    # Copy the formatting (style) of the last existing data row (224) onto
    # the new row's date cell so it keeps the same date number format /
    # border / alignment as the rest of column A, without inventing a new
    # cell style.
    $ws.Range("A224").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Range("A$r").Value = $item.Date
    $ws.Range("B$r").Value = "ECONOMICS:ROM2"
    $ws.Range("C$r").Value = $item.Value
    $ws.Range("D$r").Value = $item.Value
    $ws.Range("E$r").Value = $item.Value
    $ws.Range("F$r").Value = $item.Value
    $ws.Range("G$r").Value = 0
}
